# Remove the obsolete "RESPONSE / Invalid epoch timestamp" testcase that used
# to live in rows 5-6 (PUT with validity_in_minutes=10, key_validity_start_in_epoch=-1,
# followed by its ASSERTCONTAINS check). Deleting the two entire rows shifts the
# remaining testcase (previously rows 7-8) up to become the new rows 5-6.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("5:6").EntireRow.Delete() | Out-Null

# Leave the sheet selection where Excel would land after deleting those rows.
$ws.Activate()
$ws.Range("B10").Select() | Out-Null
